$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '304.50'
Set-TextValue 'E2' '0.83%'
Set-TextValue 'D3' '35.98'
Set-TextValue 'E3' '-3.86%'
Set-TextValue 'E4' '1.81%'
Set-TextValue 'D5' '0.07858'
Set-TextValue 'E5' '0.40%'
Set-TextValue 'D6' '2.115'
Set-TextValue 'E6' '-3.73%'
Set-TextValue 'D7' '7.949'
Set-TextValue 'E7' '-1.01%'
Set-TextValue 'B8' 'MXToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D8' '0.9212'
Set-TextValue 'E8' '1.30%'
Set-TextValue 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D9' '0.09749'
Set-TextValue 'E9' '0.91%'
Set-TextValue 'B10' 'WazirX'
Set-TextValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1856'
Set-TextValue 'E10' '-1.59%'
Set-TextValue 'B11' 'MandalaExchangeToken'
Set-TextValue 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.08682'
Set-TextValue 'E11' '1.99%'
Set-TextValue 'B12' 'BitrueCoin'
Set-TextValue 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D12' '0.03549'
Set-TextValue 'E12' '0.69%'
Set-TextValue 'B13' 'BitMartToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D13' '0.09943'
Set-TextValue 'E13' '-0.15%'
Set-TextValue 'B14' 'BitForexToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D14' '0.001433'
Set-TextValue 'E14' '-3.48%'
Set-TextValue 'B15' 'TigerCash'
Set-TextValue 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D15' '0.005701'
Set-TextValue 'E15' '0.00%'
Set-TextValue 'B16' 'LEO'
Set-TextValue 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D16' '3.467'
Set-TextValue 'E16' '0.10%'
Set-TextValue 'B17' 'GateToken'
Set-TextValue 'C17' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D17' '4.104'
Set-TextValue 'E17' '2.13%'
Set-TextValue 'D18' '2.466'
Set-TextValue 'E18' '19.16%'
Set-TextValue 'D19' '0.3423'
Set-TextValue 'E19' '-1.16%'
Set-TextValue 'D20' '5.278'
Set-TextValue 'E20' '10.79%'
Set-TextValue 'D21' '0.1302'
Set-TextValue 'E21' '0.67%'
Set-TextValue 'E22' '-0.20%'
Set-TextValue 'D23' '0.04553'
Set-TextValue 'E23' '-1.63%'
Set-TextValue 'D24' '0.005082'
Set-TextValue 'E24' '5.90%'
Set-TextValue 'D25' '0.001236'
Set-TextValue 'E25' '0.45%'
Set-TextValue 'D27' '0.0004758'
Set-TextValue 'E27' '0.07%'
Set-TextValue 'D39' '0.01841'
Set-TextValue 'E39' '5.20%'
Set-TextValue 'D40' '0.04724'
Set-TextValue 'E40' '0.35%'
Set-TextValue 'E41' '-6.64%'
Set-TextValue 'D42' '0.1403'
Set-TextValue 'E42' '0.72%'
Set-TextValue 'D43' '0.007746'
Set-TextValue 'E43' '0.94%'
Set-TextValue 'D44' '0.002218'
Set-TextValue 'E44' '-0.66%'
Set-TextValue 'D45' '0.01101'
Set-TextValue 'E45' '7.62%'
Set-TextValue 'D46' '0.00006336'
Set-TextValue 'E46' '4.93%'
Set-TextValue 'E47' '0.08%'
Set-TextValue 'D48' '0.0005811'
Set-TextValue 'E48' '0.18%'
Set-TextValue 'D49' '41.30'
Set-TextValue 'E49' '376.35%'
Set-TextValue 'D50' '0.002004'
Set-TextValue 'E50' '-25.53%'
Set-TextValue 'D51' '0.00002104'
Set-TextValue 'E51' '0.08%'
